$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Range("A7").Value = 131085613
$ws.Range("B7").Value = 57884
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "färska spår"
$ws.Range("Q7").Value = 585222
$ws.Range("R7").Value = 7060481
$ws.Range("S7").Value = 15
$ws.Range("Z7").Value = "11:52"
$ws.Range("AB7").Value = "11:52"
$ws.Range("AC7").Value = "Färska ringhack gran"
$ws.Range("AW7").Value = "Daniel Rutschman"
$ws.Range("AX7").Value = "Daniel Rutschman"
# --- Row 8 ---
$ws.Range("A8").Value = 131085805
$ws.Range("B8").Value = 79244
$ws.Range("Q8").Value = 585215
$ws.Range("R8").Value = 7060513
$ws.Range("S8").Value = 10
$ws.Range("Z8").Value = "12:01"
$ws.Range("AB8").Value = "12:01"
$ws.Range("AW8").Value = "Kim Hultgren"
$ws.Range("AX8").Value = "Kim Hultgren"
# --- Row 9 ---
$ws.Range("A9").Value = 131092646
$ws.Range("B9").Value = 79244
$ws.Range("Q9").Value = 585082
$ws.Range("R9").Value = 7060264
$ws.Range("S9").Value = 15
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
$ws.Range("AW9").Value = "Daniel Rutschman"
$ws.Range("AX9").Value = "Daniel Rutschman"
# --- Row 10 ---
$ws.Range("A10").Value = 131086958
$ws.Range("B10").Value = 79244
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value = 585165
$ws.Range("R10").Value = 7060565
$ws.Range("S10").Value = 10
$ws.Range("Z10").Value = "12:22"
$ws.Range("AB10").Value = "12:22"
$ws.Range("AC10").ClearContents()
$ws.Range("AW10").Value = "Kim Hultgren"
$ws.Range("AX10").Value = "Kim Hultgren"
# --- Row 11 ---
$ws.Range("A11").Value = 131087481
$ws.Range("B11").Value = 91829
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H11").Value = ""
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 585150
$ws.Range("R11").Value = 7060657
$ws.Range("S11").Value = 15
$ws.Range("Z11").ClearContents()
$ws.Range("AB11").ClearContents()
$ws.Range("AC11").ClearContents()
$ws.Range("AW11").Value = "Daniel Rutschman"
$ws.Range("AX11").Value = "Daniel Rutschman"
# --- Row 12 ---
$ws.Range("A12").Value = 131085696
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "färska spår"
$ws.Range("Q12").Value = 585207
$ws.Range("R12").Value = 7060471
$ws.Range("S12").Value = 10
$ws.Range("Z12").Value = "11:55"
$ws.Range("AB12").Value = "11:55"
$ws.Range("AC12").Value = "Ringhack på gran"
$ws.Range("AW12").Value = "Kim Hultgren"
$ws.Range("AX12").Value = "Kim Hultgren"
# --- Row 13 ---
$ws.Range("B13").Value = 91829
# --- Row 14 ---
$ws.Range("B14").Value = 91829
# --- Row 16 ---
$ws.Range("B16").Value = 91829
# --- Row 18 ---
$ws.Range("B18").Value = 91805
# --- Row 19 ---
$ws.Range("B19").Value = 79244
# --- Row 20 ---
$ws.Range("B20").Value = 79244
# --- Row 23 ---
$ws.Range("B23").Value = 79244
# --- Row 25 ---
$ws.Range("B25").Value = 79244
# --- Row 28 ---
$ws.Range("B28").Value = 91805
# --- Row 29 ---
$ws.Range("B29").Value = 91829
# --- Row 32 ---
$ws.Range("B32").Value = 91805
# --- Row 33 ---
$ws.Range("A33").Value = 131085569
$ws.Range("B33").Value = 79244
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("M33").ClearContents()
$ws.Range("Q33").Value = 585249
$ws.Range("R33").Value = 7060505
# --- Row 34 ---
$ws.Range("A34").Value = 131087388
$ws.Range("B34").Value = 79244
$ws.Range("Q34").Value = 585131
$ws.Range("R34").Value = 7060627
# --- Row 35 ---
$ws.Range("A35").Value = 131092590
$ws.Range("B35").Value = 79244
$ws.Range("Q35").Value = 585145
$ws.Range("R35").Value = 7060230
$ws.Range("S35").Value = 10
$ws.Range("Z35").Value = "15:20"
$ws.Range("AB35").Value = "15:20"
$ws.Range("AW35").Value = "Kim Hultgren"
$ws.Range("AX35").Value = "Kim Hultgren"
# --- Row 36 ---
$ws.Range("A36").Value = 131085126
$ws.Range("B36").Value = 57064
$ws.Range("E36").Value = 102612
$ws.Range("F36").Value = "Järpe"
$ws.Range("G36").Value = "Tetrastes bonasia"
$ws.Range("H36").Value = "(Linnaeus, 1758)"
$ws.Range("M36").Value = "färsk spillning"
$ws.Range("Q36").Value = 585219
$ws.Range("R36").Value = 7060240
$ws.Range("S36").Value = 15
$ws.Range("Z36").ClearContents()
$ws.Range("AB36").ClearContents()
$ws.Range("AW36").Value = "Daniel Rutschman"
$ws.Range("AX36").Value = "Daniel Rutschman"

Write-Host "Applied all changes"